# Adds the "getWeekenders.py" row (new row 13) to the "scripts (.py)" overview
# sheet, matching the commit "Added getWeekenders to overview, display change
# in getWeekenders".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row -----------------------------------------------------
$ws.Range("A13").Value = "streak_compuation"
$ws.Range("B13").Value = "getWeekenders.py"
$ws.Range("C13").Value = "commitsdates_per_user.json, users_reduced.json"
$ws.Range("D13").Value = "weekenders.json"
$ws.Range("E13").Value = '{"[userID]": {"WD": [numberContributions], "WE": [numberContributions]}, … }'
$ws.Range("F13").Value = '{"712": {"WD": 1627, "WE": 524}, … }'
$ws.Range("H13").Value = "creates file for number of contributions on weekdays and weekends (in the local timezone) for each user"

# Row 13 wraps onto three lines like the other description rows, so give it
# the same kind of explicit row height the sheet uses elsewhere.
$ws.Rows.Item(13).RowHeight = 45

# --- view / selection change -------------------------------------------
# The saved view no longer scrolls to A4 and the selection moves to the
# newly added example cell, F13.
$ws.Range("F13").Select()
